$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 430.9375
$ws.Range("I19").Value = 498.57144
$ws.Range("J19").Value = 378.33334
$ws.Range("K19").Value = 498.57144
$ws.Range("L19").Value = 378.33334
$ws.Range("M19").Value = -323.57144
$ws.Range("N19").Value = -728.33334
# Row 116
$ws.Range("H116").Value = 7268.9473
$ws.Range("I116").Value = 10159.167
$ws.Range("J116").Value = 2314.2856
$ws.Range("K116").Value = 10159.167
$ws.Range("L116").Value = 2314.2856
$ws.Range("M116").Value = -6717.166999999999
$ws.Range("N116").Value = -9198.285599999999
# Row 135
$ws.Range("H135").Value = 1407.8077
$ws.Range("I135").Value = 1402.68
$ws.Range("K135").Value = 12624.12
$ws.Range("M135").Value = -10089.12

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6571.8184
$ws.Range("I61").Value = 6741.905
$ws.Range("K61").Value = 6741.905
$ws.Range("M61").Value = -6529.905
# Row 74
$ws.Range("H74").Value = 1589.4203
$ws.Range("I74").Value = 1468.871
$ws.Range("J74").Value = 2657.1428
$ws.Range("K74").Value = 1468.871
$ws.Range("L74").Value = 2657.1428
$ws.Range("M74").Value = -594.8710000000001
$ws.Range("N74").Value = -4405.1428
# Row 77
$ws.Range("H77").Value = 1589.4203
$ws.Range("I77").Value = 1468.871
$ws.Range("J77").Value = 2657.1428
$ws.Range("K77").Value = 7344.355
$ws.Range("L77").Value = 13285.714
$ws.Range("M77").Value = -2976.355
$ws.Range("N77").Value = -22021.714
# Row 110
$ws.Range("H110").Value = 672.0345
$ws.Range("I110").Value = 587.5263
$ws.Range("J110").Value = 832.6
$ws.Range("K110").Value = 587.5263
$ws.Range("L110").Value = 832.6
$ws.Range("M110").Value = 1457.4737
$ws.Range("N110").Value = -4922.6
# Row 132
$ws.Range("H132").Value = 3570.7273
$ws.Range("I132").Value = 1679.1904
$ws.Range("J132").Value = 6880.9165
$ws.Range("K132").Value = 5037.5712
$ws.Range("L132").Value = 20642.7495
$ws.Range("M132").Value = -2507.5712
$ws.Range("N132").Value = -25702.7495
# Row 136
$ws.Range("H136").Value = 6571.8184
$ws.Range("I136").Value = 6741.905
$ws.Range("K136").Value = 20225.715
$ws.Range("M136").Value = -17675.715

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 990.24243
$ws.Range("I94").Value = 649.2143
$ws.Range("J94").Value = 2900
$ws.Range("K94").Value = 649.2143
$ws.Range("L94").Value = 2900
$ws.Range("M94").Value = -198.2143
$ws.Range("N94").Value = -3802
# Row 134
$ws.Range("H134").Value = 7505.364
$ws.Range("I134").Value = 11235.75
$ws.Range("J134").Value = 3028.9
$ws.Range("K134").Value = 33707.25
$ws.Range("L134").Value = 9086.700000000001
$ws.Range("M134").Value = -31172.25
$ws.Range("N134").Value = -14156.7

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 369.72726
$ws.Range("I107").Value = 816.75
$ws.Range("K107").Value = 816.75
$ws.Range("M107").Value = 1103.25
# Row 122
$ws.Range("H122").Value = 952
$ws.Range("I122").Value = 952
$ws.Range("K122").Value = 2856
$ws.Range("M122").Value = -406
# Row 134
$ws.Range("H134").Value = 2870.8064
$ws.Range("I134").Value = 3679.3157
$ws.Range("K134").Value = 11037.9471
$ws.Range("M134").Value = -8502.947100000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 5816.5
$ws.Range("J55").Value = 6579.8
$ws.Range("L55").Value = 19739.4
$ws.Range("N55").Value = -20093.4
# Row 131
$ws.Range("H131").Value = 2858039.5
$ws.Range("I131").Value = 7143577.5
$ws.Range("J131").Value = 1014.4286
$ws.Range("K131").Value = 21430732.5
$ws.Range("L131").Value = 3043.2858
$ws.Range("M131").Value = -21425692.5
$ws.Range("N131").Value = -13123.2858

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Range("H47").Value = 4031
$ws.Range("J47").Value = 4031
$ws.Range("L47").Value = 4031
$ws.Range("N47").Value = -5167
# Row 80
$ws.Range("H80").Value = 2309.9092
$ws.Range("I80").Value = 1321.8
$ws.Range("J80").Value = 3133.3333
$ws.Range("K80").Value = 1321.8
$ws.Range("L80").Value = 3133.3333
$ws.Range("M80").Value = -323.8
$ws.Range("N80").Value = -5129.3333
# Row 83
$ws.Range("H83").Value = 2309.9092
$ws.Range("I83").Value = 1321.8
$ws.Range("J83").Value = 3133.3333
$ws.Range("K83").Value = 6609
$ws.Range("L83").Value = 15666.6665
$ws.Range("M83").Value = -1617
$ws.Range("N83").Value = -25650.6665
# Row 132
$ws.Range("H132").Value = 3116.56
$ws.Range("I132").Value = 2728.4
$ws.Range("J132").Value = 3698.8
$ws.Range("K132").Value = 8185.200000000001
$ws.Range("L132").Value = 11096.4
$ws.Range("M132").Value = -5655.200000000001
$ws.Range("N132").Value = -16156.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 125001760
$ws.Range("I40").Value = 125001760
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 125001760
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -125001624
$ws.Range("N40").ClearContents()
# Row 46
$ws.Range("H46").Value = 11905613
$ws.Range("I46").Value = 18519176
$ws.Range("J46").Value = 1199
$ws.Range("K46").Value = 18519176
$ws.Range("L46").Value = 1199
$ws.Range("M46").Value = -18518988
$ws.Range("N46").Value = -1575
# Row 61
$ws.Range("H61").Value = 2192.4167
$ws.Range("I61").Value = 1882.6
$ws.Range("J61").Value = 2708.7778
$ws.Range("K61").Value = 1882.6
$ws.Range("L61").Value = 2708.7778
$ws.Range("M61").Value = -1680.6
$ws.Range("N61").Value = -3112.7778
# Row 100
$ws.Range("H100").Value = 2647.818
$ws.Range("I100").Value = 2200.3333
$ws.Range("K100").Value = 2200.3333
$ws.Range("M100").Value = -1659.3333
# Row 113
$ws.Range("H113").Value = 2192.4167
$ws.Range("I113").Value = 1882.6
$ws.Range("J113").Value = 2708.7778
$ws.Range("K113").Value = 1882.6
$ws.Range("L113").Value = 2708.7778
$ws.Range("M113").Value = 287.4000000000001
$ws.Range("N113").Value = -7048.7778
# Row 132
$ws.Range("H132").Value = 16056283
$ws.Range("I132").Value = 19704252
$ws.Range("K132").Value = 59112756
$ws.Range("M132").Value = -59110226
# Row 136
$ws.Range("H136").Value = 5691.0225
$ws.Range("I136").Value = 5177.528
$ws.Range("K136").Value = 15532.584
$ws.Range("M136").Value = -12982.584

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1519.4166
$ws.Range("I122").Value = 1603.3
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 4809.9
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -2359.9
$ws.Range("N122").Value = -8200
